# Updated project to parse new input excel file.
# - Rename unit-bearing headers to remove parentheses (e.g. "Electricity (KWH)" -> "Electricity_KWH")
# - Replace the text Month column values with numeric month-of-year values
# - Add a new numeric "Year" column (G) derived from each row's month text
#   (Nov/Dec entries belong to 2015, all other months belong to 2016)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames -------------------------------------------------------
$ws.Range("B1").Value = "Electricity_KWH"
$ws.Range("C1").Value = "Cold_Water_m3"
$ws.Range("D1").Value = "Hot_Water_m3"
$ws.Range("E1").Value = "Heat_MWH"
$ws.Range("G1").Value = "Year"

# --- Month name -> (month number, year) lookup ----------------------------
$monthMap = @{
    "January"   = @(1, 2016)
    "February"  = @(2, 2016)
    "March"     = @(3, 2016)
    "April"     = @(4, 2016)
    "May"       = @(5, 2016)
    "June"      = @(6, 2016)
    "July"      = @(7, 2016)
    "August"    = @(8, 2016)
    "September" = @(9, 2016)
    "October"   = @(10, 2016)
    "November"  = @(11, 2015)
    "December"  = @(12, 2015)
}

# --- Walk every data row and replace Month text with numeric Month + Year -
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $monthText = $ws.Cells.Item($r, 6).Text
    $pair = $monthMap[$monthText]
    $ws.Cells.Item($r, 6).Value = $pair[0]
    $ws.Cells.Item($r, 7).Value = $pair[1]
}
